$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume number, week-covering dates) ---
$ws.Range("A8").Value = "Volume 31   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/8/2024  Through  1/14/2024"

# --- Cells whose type/style flips between numeric and the text placeholders ("0" / "***.*") ---
$ws.Range("J14").Copy($ws.Range("C14"))
$ws.Range("J14").Copy($ws.Range("G14"))
$ws.Range("K14").Copy($ws.Range("H14"))
$ws.Range("J14").Copy($ws.Range("C15"))
$ws.Range("F14").Copy($ws.Range("D15"))
$ws.Range("D15").Value2 = 1
$ws.Range("L14").Copy($ws.Range("E15"))
$ws.Range("E15").Value2 = -100
$ws.Range("F14").Copy($ws.Range("J15"))
$ws.Range("J15").Value2 = 1
$ws.Range("L14").Copy($ws.Range("K15"))
$ws.Range("K15").Value2 = 0
$ws.Range("F14").Copy($ws.Range("C20"))
$ws.Range("C20").Value2 = 3
$ws.Range("F14").Copy($ws.Range("I20"))
$ws.Range("I20").Value2 = 3
$ws.Range("J14").Copy($ws.Range("C22"))
$ws.Range("J14").Copy($ws.Range("D22"))
$ws.Range("K14").Copy($ws.Range("E22"))
$ws.Range("J14").Copy($ws.Range("F23"))
$ws.Range("J14").Copy($ws.Range("C26"))
$ws.Range("F14").Copy($ws.Range("D26"))
$ws.Range("D26").Value2 = 4
$ws.Range("L14").Copy($ws.Range("E26"))
$ws.Range("E26").Value2 = -100
$ws.Range("F14").Copy($ws.Range("J26"))
$ws.Range("J26").Value2 = 4
$ws.Range("L14").Copy($ws.Range("K26"))
$ws.Range("K26").Value2 = -75
$ws.Range("L14").Copy($ws.Range("M28"))
$ws.Range("M28").Value2 = -100
$ws.Range("L14").Copy($ws.Range("M29"))
$ws.Range("M29").Value2 = -100

# --- Plain numeric value updates (style unchanged) ---
$ws.Range("N14").Value2 = -50
$ws.Range("F15").Value2 = 2
$ws.Range("G15").Value2 = 4
$ws.Range("H15").Value2 = -50
$ws.Range("C16").Value2 = 11
$ws.Range("D16").Value2 = 6
$ws.Range("E16").Value2 = 83.333333333333
$ws.Range("F16").Value2 = 37
$ws.Range("G16").Value2 = 28
$ws.Range("H16").Value2 = 32.142857142857
$ws.Range("I16").Value2 = 18
$ws.Range("J16").Value2 = 13
$ws.Range("K16").Value2 = 38.461538461538
$ws.Range("L16").Value2 = -10
$ws.Range("M16").Value2 = -14.285714285714
$ws.Range("N16").Value2 = -81.25
$ws.Range("C17").Value2 = 13
$ws.Range("D17").Value2 = 8
$ws.Range("E17").Value2 = 62.5
$ws.Range("F17").Value2 = 54
$ws.Range("G17").Value2 = 36
$ws.Range("H17").Value2 = 50
$ws.Range("I17").Value2 = 27
$ws.Range("J17").Value2 = 17
$ws.Range("K17").Value2 = 58.823529411764
$ws.Range("L17").Value2 = -10
$ws.Range("M17").Value2 = 170
$ws.Range("N17").Value2 = -35.714285714285
$ws.Range("C18").Value2 = 4
$ws.Range("D18").Value2 = 3
$ws.Range("E18").Value2 = 33.333333333333
$ws.Range("F18").Value2 = 22
$ws.Range("G18").Value2 = 15
$ws.Range("H18").Value2 = 46.666666666666
$ws.Range("I18").Value2 = 7
$ws.Range("J18").Value2 = 8
$ws.Range("K18").Value2 = -12.5
$ws.Range("L18").Value2 = 0
$ws.Range("M18").Value2 = -50
$ws.Range("N18").Value2 = -91.566265060241
$ws.Range("C19").Value2 = 13
$ws.Range("D19").Value2 = 7
$ws.Range("E19").Value2 = 85.714285714285
$ws.Range("F19").Value2 = 45
$ws.Range("G19").Value2 = 53
$ws.Range("H19").Value2 = -15.094339622641
$ws.Range("I19").Value2 = 25
$ws.Range("J19").Value2 = 19
$ws.Range("K19").Value2 = 31.578947368421
$ws.Range("L19").Value2 = 8.695652173913
$ws.Range("M19").Value2 = 316.666666666667
$ws.Range("N19").Value2 = -26.470588235294
$ws.Range("D20").Value2 = 5
$ws.Range("E20").Value2 = -40
$ws.Range("F20").Value2 = 15
$ws.Range("G20").Value2 = 21
$ws.Range("H20").Value2 = -28.571428571428
$ws.Range("J20").Value2 = 12
$ws.Range("K20").Value2 = -75
$ws.Range("L20").Value2 = -75
$ws.Range("M20").Value2 = -25
$ws.Range("N20").Value2 = -91.891891891891
$ws.Range("C21").Value2 = 44
$ws.Range("D21").Value2 = 30
$ws.Range("E21").Value2 = 46.666666666666
$ws.Range("F21").Value2 = 176
$ws.Range("G21").Value2 = 157
$ws.Range("H21").Value2 = 12.101910828025
$ws.Range("I21").Value2 = 82
$ws.Range("J21").Value2 = 70
$ws.Range("K21").Value2 = 17.142857142857
$ws.Range("L21").Value2 = -12.765957446808
$ws.Range("M21").Value2 = 49.090909090909
$ws.Range("N21").Value2 = -72.108843537415
$ws.Range("H23").Value2 = -100
$ws.Range("J23").Value2 = 2
$ws.Range("C24").Value2 = 20
$ws.Range("D24").Value2 = 22
$ws.Range("E24").Value2 = -9.090909090909
$ws.Range("F24").Value2 = 67
$ws.Range("H24").Value2 = 1.515151515151
$ws.Range("I24").Value2 = 30
$ws.Range("J24").Value2 = 36
$ws.Range("K24").Value2 = -16.666666666666
$ws.Range("L24").Value2 = -21.052631578947
$ws.Range("M24").Value2 = -6.25
$ws.Range("C25").Value2 = 29
$ws.Range("D25").Value2 = 23
$ws.Range("E25").Value2 = 26.086956521739
$ws.Range("F25").Value2 = 98
$ws.Range("G25").Value2 = 82
$ws.Range("H25").Value2 = 19.512195121951
$ws.Range("I25").Value2 = 47
$ws.Range("J25").Value2 = 45
$ws.Range("K25").Value2 = 4.444444444444
$ws.Range("L25").Value2 = 113.636363636364
$ws.Range("M25").Value2 = 20.512820512820
$ws.Range("F26").Value2 = 5
$ws.Range("G26").Value2 = 8
$ws.Range("H26").Value2 = -37.5
$ws.Range("C27").Value2 = 4
$ws.Range("D27").Value2 = 2
$ws.Range("F27").Value2 = 8
$ws.Range("G27").Value2 = 4
$ws.Range("H27").Value2 = 100
$ws.Range("I27").Value2 = 6
$ws.Range("J27").Value2 = 3
$ws.Range("F28").Value2 = 1
$ws.Range("G28").Value2 = 1
$ws.Range("H28").Value2 = 0
$ws.Range("F29").Value2 = 1
$ws.Range("G29").Value2 = 1
$ws.Range("H29").Value2 = 0
$ws.Range("J42").Value2 = 279
$ws.Range("K42").Value2 = -34.660421545667
$ws.Range("L42").Value2 = -32.281553398058
$ws.Range("M42").Value2 = -68.114285714285
$ws.Range("N42").Value2 = -77.481840193704
$ws.Range("J43").Value2 = 2193
$ws.Range("K43").Value2 = -20.254545454545
$ws.Range("L43").Value2 = -38.120767494356
$ws.Range("M43").Value2 = -70.818363273453
$ws.Range("N43").Value2 = -74.822043628013
